$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.619.94'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '1.665.44'
$ws.Range('E3').Value = '  -3.91%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '215.56'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('E6').Value = '  -2.51%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.20'
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('D11').Value = '0.0878'
$ws.Range('D12').Value = '1.901.00'
$ws.Range('E12').Value = '  -4.01%  '
$ws.Range('D13').Value = '1.702.76'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('D15').Value = '0.571'
$ws.Range('E15').Value = '  +1.25%  '
$ws.Range('D16').Value = '66.31'
$ws.Range('E16').Value = '  -1.99%  '
$ws.Range('D17').Value = '27.601.33'
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('D18').Value = '242.08'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '0.0₃0732'
$ws.Range('E19').Value = '  -3.29%  '
$ws.Range('D20').Value = '7.69'
$ws.Range('E20').Value = '  -3.83%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '4.52'
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').Value = '9.38'
$ws.Range('E23').Value = '  -3.43%  '
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('D25').Value = '146.47'
$ws.Range('E25').Value = '  -2.24%  '
$ws.Range('D26').Value = '7.28'
$ws.Range('E26').Value = '  -3.41%  '
$ws.Range('D27').Value = '16.38'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('D30').Value = '1.22'
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('D33').Value = '1.457.85'
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('E35').Value = '  -3.74%  '
$ws.Range('E36').Value = '  -3.60%  '
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('D38').Value = '0.576'
$ws.Range('E38').Value = '  -4.84%  '
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('D40').Value = '69.94'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = '2.23'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('D44').Value = '5.42'
$ws.Range('E44').Value = '  -5.21%  '
$ws.Range('D45').Value = '0.794'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('D46').Value = '1.808.86'
$ws.Range('E46').Value = '  -3.97%  '
$ws.Range('D47').Value = '1.73'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '88.92'
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('E49').Value = '  -6.08%  '
$ws.Range('E50').Value = '  -1.82%  '
$ws.Range('E51').Value = '  -4.08%  '
